# Remmoving past portfolio data
# Shifts the consumption forecast window forward by 4 days:
#   - Timestamps (column B) are each advanced by 4 days.
#   - Forecasted consumption values (column A) are replaced with the
#     new forecast data for the new window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for rows 2..97 (A column)
$newConsumption = @(
    5210,5170,5130,5090,5050,5030,5010,5000,5000,5000,
    5000,5010,5030,5050,5080,5140,5200,5290,5400,5530,
    5730,5910,6090,6270,6470,6640,6790,6920,7040,7090,
    7110,7110,7090,7040,6980,6900,6820,6740,6660,6590,
    6520,6470,6420,6370,6320,6300,6290,6280,6280,6280,
    6280,6280,6290,6300,6310,6330,6390,6450,6500,6560,
    6620,6700,6780,6860,6960,7050,7150,7250,7380,7470,
    7570,7660,7750,7810,7820,7810,7770,7700,7640,7570,
    7470,7380,7280,7150,6970,6820,6660,6530,6410,6280,
    6170,6040,5950,5900,5860,5810
)

$firstRow = 2
$lastRow = 97

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $idx = $row - $firstRow

    # Column A: new forecast value
    $ws.Cells.Item($row, 1).Value = $newConsumption[$idx]

    # Column B: shift the existing timestamp forward by 4 days
    $tsCell = $ws.Cells.Item($row, 2)
    $oldTs = $tsCell.Value()
    $tsCell.Value = $oldTs.AddDays(4)
}
